$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 65
$ws.Range("D3").Value = 99
$ws.Range("H3").Value = 96
$ws.Range("C6").Value = 351
$ws.Range("D6").Value = 308
$ws.Range("F6").Value = 394
$ws.Range("H6").Value = 325
$ws.Range("I6").Value = 384
$ws.Range("C7").Value = 473
$ws.Range("D7").Value = 481
$ws.Range("F7").Value = 559
$ws.Range("H7").Value = 516
$ws.Range("I7").Value = 636

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("H6").Value = 31
$ws.Range("H7").Value = 40

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F2").Value = 3
$ws.Range("F6").Value = 9

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 7

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("C4").Value = 2
$ws.Range("F5").Value = 9
$ws.Range("D7").Value = 7
$ws.Range("H8").Value = 38
$ws.Range("H32").Value = 40
$ws.Range("C35").Value = 8
$ws.Range("F52").Value = 3
$ws.Range("I52").Value = 4
$ws.Range("D53").Value = 59
$ws.Range("C80").Value = 7
$ws.Range("C98").Value = 473
$ws.Range("D98").Value = 481
$ws.Range("F98").Value = 559
$ws.Range("H98").Value = 516
$ws.Range("I98").Value = 636

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D3").Value = 16
$ws.Range("D7").Value = 59

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 8

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 7

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 2

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("F5").Value = 3
$ws.Range("I5").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("I6").Value = 4

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("H3").Value = 4
$ws.Range("H6").Value = 38
